# Add the "SRC_By_Day" worksheet (tab-delimited SRC-by-day export) after the
# existing "SupplyDemand" sheet, and populate it with the forge data.

$wb = $excel.ActiveWorkbook

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "SRC_By_Day"

# ---------------------------------------------------------------------
# Row 1 - phase headers (sparse, spans visually over the day columns)
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "FwdStation"
$ws.Range("H1").Value = "PH I"
$ws.Range("L1").Value = "PH IIa"
$ws.Range("P1").Value = "PH IIb"

# ---------------------------------------------------------------------
# Row 2 - column headers
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "SRC"
$ws.Range("B2").Value = "Title"
$ws.Range("C2").Value = "Strength"
$ws.Range("D2").Value = "Branc Code"
$ws.Range("E2").Value = "Branch Label"
$ws.Range("F2").Value = "Service"
$ws.Range("G2").Value = "Day 0001`nTP 1"
$ws.Range("H2").Value = "Day 0009`nTP 2"
$ws.Range("I2").Value = "Day 0017`nTP 3"
$ws.Range("J2").Value = "Day 0025`nTP 4"
$ws.Range("K2").Value = "Day 0033`nTP 5"
$ws.Range("L2").Value = "Day 0041`nTP 6"
$ws.Range("M2").Value = "Day 0049`nTP 7"
$ws.Range("N2").Value = "Day 0057`nTP 8"
$ws.Range("O2").Value = "Day 0065`nTP 9"
$ws.Range("P2").Value = "Day 0073`nTP 10"
$ws.Range("Q2").Value = "Day 0081`nTP 11"
$ws.Range("R2").Value = "Day 0089`nTP 12"

# day-columns on the header row wrap their two-line labels
$hdr = $ws.Range("G2:R2")
$hdr.WrapText = $true
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 10
$ws.Rows.Item(2).RowHeight = 23.85

# ---------------------------------------------------------------------
# Row 3 - SRC 01205K000 / ASLT BN
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "01205K000"
$ws.Range("B3").Value = "Ponies"
$ws.Range("C3").Value = 381
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Avn"
$ws.Range("F3").Value = "Army"
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 4
$ws.Range("R3").Value = 7

# ---------------------------------------------------------------------
# Row 4 - SRC 01225K000 / GSAB
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "01225K000"
$ws.Range("B4").Value = "Buffalos"
$ws.Range("C4").Value = 590
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Avn"
$ws.Range("F4").Value = "Army"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 8
$ws.Range("P4").Value = 8
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 7

# ---------------------------------------------------------------------
# Row 6 - totals / TP index row
# ---------------------------------------------------------------------
$ws.Range("B6").Value = "Total Army Strength"
$ws.Range("C6").Value = 362377
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 7
$ws.Range("N6").Value = 8
$ws.Range("O6").Value = 9
$ws.Range("P6").Value = 10
$ws.Range("Q6").Value = 11
$ws.Range("R6").Value = 12

# ---------------------------------------------------------------------
# column widths (approximate character widths from the source export)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.77
$ws.Columns.Item(2).ColumnWidth = 16.72
$ws.Columns.Item(3).ColumnWidth = 7.55
$ws.Columns.Item(4).ColumnWidth = 10.2
$ws.Columns.Item(5).ColumnWidth = 11.3
$ws.Columns.Item(6).ColumnWidth = 6.85
for ($c = 7; $c -le 15; $c++) {
  $ws.Columns.Item($c).ColumnWidth = 12.13
}
$ws.Columns.Item(16).ColumnWidth = 13.1
$ws.Columns.Item(17).ColumnWidth = 12.98
$ws.Columns.Item(18).ColumnWidth = 13.1

$ws.Range("A1").Select()
